$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "37.107.72"
$cell.ClearFormats()
$ws.Range("E2").Value = "  -1.60%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.992.62"
$cell.ClearFormats()
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  +0.00%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "242.56"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -6.27%  "
$ws.Range("E6").Value = "  -3.86%  "
$ws.Range("E7").Value = "  +0.06%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "54.58"
$cell.ClearFormats()
$ws.Range("E8").Value = "  -5.57%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.371"
$cell.ClearFormats()
$ws.Range("E9").Value = "  -4.53%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "58.66"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +2.35%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0751"
$cell.ClearFormats()
$ws.Range("E11").Value = "  -6.19%  "
$ws.Range("E12").Value = "  -4.32%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "2.284.14"
$cell.ClearFormats()
$ws.Range("E13").Value = "  -2.57%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "13.98"
$cell.ClearFormats()
$ws.Range("E14").Value = "  -5.68%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "21.01"
$cell.ClearFormats()
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("E16").Value = "  -8.29%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "5.05"
$cell.ClearFormats()
$ws.Range("E17").Value = "  -6.21%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "1.992.17"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -2.62%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "37.004.61"
$cell.ClearFormats()
$ws.Range("E19").Value = "  -1.48%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "68.19"
$cell.ClearFormats()
$ws.Range("E20").Value = "  -2.91%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0810"
$cell.ClearFormats()
$ws.Range("E21").Value = "  -5.58%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "226.84"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("E23").Value = "  -5.34%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  -9.44%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.36"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +0.07%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "161.68"
$cell.ClearFormats()
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("E28").Value = "  -5.67%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "19.09"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -4.83%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.122"
$cell.ClearFormats()
$ws.Range("E30").Value = "  -11.27%  "
$ws.Range("E31").Value = "  -5.26%  "
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("E33").Value = "  -7.44%  "
$ws.Range("E34").Value = "  -8.17%  "
$ws.Range("E35").Value = "  -6.34%  "
$ws.Range("E36").Value = "  -6.91%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  -1.67%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "3.31"
$cell.ClearFormats()
$ws.Range("E39").Value = "  -4.28%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "5.24"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -2.82%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "3.03"
$cell.ClearFormats()
$ws.Range("E41").Value = "  -0.13%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.429.17"
$cell.ClearFormats()
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.12"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -6.11%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.0203"
$cell.ClearFormats()
$ws.Range("E44").Value = "  -6.86%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0887"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -8.53%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "88.30"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -3.57%  "
$ws.Range("E47").Value = "  -5.32%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()
$ws.Range("E48").Value = "  -5.12%  "
$ws.Range("E49").Value = "  -0.14%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "6.63"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -10.92%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.176.13"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -2.62%  "
